# Generate Report for Handoff
# Replace the old GUID-based file name (20cbf5b9-d038-42a1-9a4e-0866393c5744)
# with the new one (a52f26a9-04f3-4155-93ce-0b26dd45f53e) across all sheets,
# and refresh the associated handoff timestamps.

$wb = $excel.ActiveWorkbook

$oldGuid = "20cbf5b9-d038-42a1-9a4e-0866393c5744"
$newGuid = "a52f26a9-04f3-4155-93ce-0b26dd45f53e"

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Hyperlinks.Item(1).TextToDisplay = "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = "2017-02-09 09:08:23"

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = "$newGuid.md"
$wsZh.Hyperlinks.Item(1).TextToDisplay = "$newGuid.md"
$wsZh.Range("G2").Value = "$newGuid.8618f1a63b792b4efe4c19a2fd11272ef5d05485.zh-cn.xlf"
$wsZh.Range("H2").Value = "2017-02-09 09:08:02"

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = "$newGuid.md"
$wsDe.Hyperlinks.Item(1).TextToDisplay = "$newGuid.md"
$wsDe.Range("G2").Value = "$newGuid.8618f1a63b792b4efe4c19a2fd11272ef5d05485.de-de.xlf"
$wsDe.Range("H2").Value = "2017-02-09 09:08:23"
